$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3783.875
$ws.Range("I19").Value = 4499
$ws.Range("J19").Value = 3681.7144
$ws.Range("K19").Value = 4499
$ws.Range("L19").Value = 3681.7144
$ws.Range("M19").Value = -4324
$ws.Range("N19").Value = -4031.7144
$ws.Range("H112").Value = 4880.8423
$ws.Range("J112").Value = 5076.5
$ws.Range("L112").Value = 15229.5
$ws.Range("N112").Value = -17445.5
$ws.Range("H116").Value = 4795.0527
$ws.Range("I116").Value = 4324.067
$ws.Range("J116").Value = 6561.25
$ws.Range("K116").Value = 4324.067
$ws.Range("L116").Value = 6561.25
$ws.Range("M116").Value = -882.067
$ws.Range("N116").Value = -13445.25
$ws.Range("H132").Value = 1967.3334
$ws.Range("I132").Value = 2023.1285
$ws.Range("J132").Value = 1479.125
$ws.Range("K132").Value = 6069.3855
$ws.Range("L132").Value = 4437.375
$ws.Range("M132").Value = -3539.3855
$ws.Range("N132").Value = -9497.375
$ws.Range("H137").Value = 2635.875
$ws.Range("I137").Value = 2417.4
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 7252.200000000001
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -4702.200000000001
$ws.Range("N137").Value = -14100
$ws.Range("H138").Value = 2859.4595
$ws.Range("I138").Value = 1296.075
$ws.Range("J138").Value = 4698.7354
$ws.Range("K138").Value = 3888.225
$ws.Range("L138").Value = 14096.2062
$ws.Range("M138").Value = 1251.775
$ws.Range("N138").Value = -24376.2062
$ws.Range("H141").Value = 755.1875
$ws.Range("I141").Value = 759.13336
$ws.Range("J141").Value = 696
$ws.Range("K141").Value = 2277.40008
$ws.Range("L141").Value = 2088
$ws.Range("M141").Value = 2902.59992
$ws.Range("N141").Value = -12448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 19835
$ws.Range("J46").Value = 19802.2
$ws.Range("L46").Value = 19802.2
$ws.Range("N46").Value = -20440.2
$ws.Range("H61").Value = 7608.4375
$ws.Range("I61").Value = 4749.3076
$ws.Range("J61").Value = 19998
$ws.Range("K61").Value = 4749.3076
$ws.Range("L61").Value = 19998
$ws.Range("M61").Value = -4537.3076
$ws.Range("N61").Value = -20422
$ws.Range("H74").Value = 1848.6666
$ws.Range("I74").Value = 1537.1538
$ws.Range("K74").Value = 1537.1538
$ws.Range("M74").Value = -663.1538
$ws.Range("H77").Value = 1848.6666
$ws.Range("I77").Value = 1537.1538
$ws.Range("K77").Value = 7685.769
$ws.Range("M77").Value = -3317.769
$ws.Range("H97").Value = 441.03125
$ws.Range("I97").Value = 465.92307
$ws.Range("J97").Value = 333.16666
$ws.Range("K97").Value = 465.92307
$ws.Range("L97").Value = 333.16666
$ws.Range("M97").Value = 30.07693
$ws.Range("N97").Value = -1325.16666
$ws.Range("H110").Value = 2111.4167
$ws.Range("I110").Value = 1961.6666
$ws.Range("J110").Value = 2261.1667
$ws.Range("K110").Value = 1961.6666
$ws.Range("L110").Value = 2261.1667
$ws.Range("M110").Value = 83.33339999999998
$ws.Range("N110").Value = -6351.1667
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 7608.4375
$ws.Range("I136").Value = 4749.3076
$ws.Range("J136").Value = 19998
$ws.Range("K136").Value = 14247.9228
$ws.Range("L136").Value = 59994
$ws.Range("M136").Value = -11697.9228
$ws.Range("N136").Value = -65094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4798.25
$ws.Range("I99").Value = 4114.5
$ws.Range("J99").Value = 6849.5
$ws.Range("K99").Value = 4114.5
$ws.Range("L99").Value = 6849.5
$ws.Range("M99").Value = -2616.5
$ws.Range("N99").Value = -9845.5
$ws.Range("H107").Value = 2902.0625
$ws.Range("I107").Value = 2292.889
$ws.Range("J107").Value = 3685.2856
$ws.Range("K107").Value = 2292.889
$ws.Range("L107").Value = 3685.2856
$ws.Range("M107").Value = -372.8890000000001
$ws.Range("N107").Value = -7525.2856
$ws.Range("H134").Value = 2138.74
$ws.Range("I134").Value = 2107.6
$ws.Range("J134").Value = 2419
$ws.Range("K134").Value = 6322.799999999999
$ws.Range("L134").Value = 7257
$ws.Range("M134").Value = -3787.799999999999
$ws.Range("N134").Value = -12327

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5413.9287
$ws.Range("I31").Value = 4425.6553
$ws.Range("J31").Value = 7618.5386
$ws.Range("K31").Value = 4425.6553
$ws.Range("L31").Value = 7618.5386
$ws.Range("M31").Value = -4130.6553
$ws.Range("N31").Value = -8208.5386
$ws.Range("H34").Value = 5413.9287
$ws.Range("I34").Value = 4425.6553
$ws.Range("J34").Value = 7618.5386
$ws.Range("K34").Value = 4425.6553
$ws.Range("L34").Value = 7618.5386
$ws.Range("M34").Value = -4223.6553
$ws.Range("N34").Value = -8022.5386
$ws.Range("H58").Value = 3664.6
$ws.Range("I58").Value = 3612.9375
$ws.Range("J58").Value = 3756.4443
$ws.Range("K58").Value = 3612.9375
$ws.Range("L58").Value = 3756.4443
$ws.Range("M58").Value = -3409.9375
$ws.Range("N58").Value = -4162.4443
$ws.Range("H99").Value = 7605.8887
$ws.Range("I99").Value = 7733.8335
$ws.Range("J99").Value = 7350
$ws.Range("K99").Value = 7733.8335
$ws.Range("L99").Value = 7350
$ws.Range("M99").Value = -6235.8335
$ws.Range("N99").Value = -10346
$ws.Range("H107").Value = 591.38464
$ws.Range("I107").Value = 502
$ws.Range("J107").Value = 1083
$ws.Range("K107").Value = 502
$ws.Range("L107").Value = 1083
$ws.Range("M107").Value = 1418
$ws.Range("N107").Value = -4923
$ws.Range("H122").Value = 252746
$ws.Range("I122").Value = 275586.72
$ws.Range("J122").Value = 1498
$ws.Range("K122").Value = 826760.1599999999
$ws.Range("L122").Value = 4494
$ws.Range("M122").Value = -824310.1599999999
$ws.Range("N122").Value = -9394
$ws.Range("H126").Value = 7605.8887
$ws.Range("I126").Value = 7733.8335
$ws.Range("J126").Value = 7350
$ws.Range("K126").Value = 23201.5005
$ws.Range("L126").Value = 22050
$ws.Range("M126").Value = -20731.5005
$ws.Range("N126").Value = -26990
$ws.Range("H132").Value = 1191.8392
$ws.Range("I132").Value = 1071.5491
$ws.Range("J132").Value = 2418.8
$ws.Range("K132").Value = 3214.6473
$ws.Range("L132").Value = 7256.400000000001
$ws.Range("M132").Value = -684.6472999999996
$ws.Range("N132").Value = -12316.4
$ws.Range("H134").Value = 1526.9272
$ws.Range("I134").Value = 1492.7646
$ws.Range("J134").Value = 1962.5
$ws.Range("K134").Value = 4478.293799999999
$ws.Range("L134").Value = 5887.5
$ws.Range("M134").Value = -1943.293799999999
$ws.Range("N134").Value = -10957.5
$ws.Range("H136").Value = 3664.6
$ws.Range("I136").Value = 3612.9375
$ws.Range("J136").Value = 3756.4443
$ws.Range("K136").Value = 10838.8125
$ws.Range("L136").Value = 11269.3329
$ws.Range("M136").Value = -8288.8125
$ws.Range("N136").Value = -16369.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 573.3333
$ws.Range("I75").Value = 360
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 1080
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = -82
$ws.Range("N75").Value = -4996
$ws.Range("H78").Value = 573.3333
$ws.Range("I78").Value = 360
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 3240
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = 1752
$ws.Range("N78").Value = -18984
$ws.Range("H113").Value = 2001
$ws.Range("J113").Value = 2271.3572
$ws.Range("L113").Value = 6814.071599999999
$ws.Range("N113").Value = -11154.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4462.2
$ws.Range("I102").Value = 4937
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 4937
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -3315
$ws.Range("N102").Value = -6994
$ws.Range("H113").Value = 5281.7393
$ws.Range("I113").Value = 5746.5293
$ws.Range("J113").Value = 3964.8333
$ws.Range("K113").Value = 5746.5293
$ws.Range("L113").Value = 3964.8333
$ws.Range("M113").Value = -3576.5293
$ws.Range("N113").Value = -8304.8333
$ws.Range("H132").Value = 6758.154
$ws.Range("I132").Value = 7883.4
$ws.Range("K132").Value = 23650.2
$ws.Range("M132").Value = -21120.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1470.3334
$ws.Range("I7").Value = 1464.4
$ws.Range("K7").Value = 1464.4
$ws.Range("M7").Value = -1352.4
$ws.Range("H126").Value = 1470.3334
$ws.Range("I126").Value = 1464.4
$ws.Range("K126").Value = 4393.200000000001
$ws.Range("M126").Value = -1923.200000000001
$ws.Range("H132").Value = 37825
$ws.Range("I132").Value = 50550
$ws.Range("K132").Value = 151650
$ws.Range("M132").Value = -149120
$ws.Range("H136").Value = 1684.1708
$ws.Range("I136").Value = 1691.275
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 5073.825000000001
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -2523.825000000001
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1039.6046
$ws.Range("I14").Value = 1039.6046
$ws.Range("K14").Value = 1039.6046
$ws.Range("M14").Value = -871.6045999999999
$ws.Range("H107").Value = 4500.8486
$ws.Range("I107").Value = 4045.3044
$ws.Range("J107").Value = 5548.6
$ws.Range("K107").Value = 12135.9132
$ws.Range("L107").Value = 16645.8
$ws.Range("M107").Value = -10215.9132
$ws.Range("N107").Value = -20485.8
$ws.Range("H122").Value = 2933.9167
$ws.Range("I122").Value = 2870.5881
$ws.Range("J122").Value = 3087.7144
$ws.Range("K122").Value = 8611.764299999999
$ws.Range("L122").Value = 9263.143199999999
$ws.Range("M122").Value = -6161.764299999999
$ws.Range("N122").Value = -14163.1432
